# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the a9c5f272-806f-4fa5-82ad-bcda24fde583 entry (row 3) across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for a9c5f272 (row 3, col G)
$wsOverview.Range("G3").Value = "2016-08-17 04:41:37"

# zh-cn sheet: "Correspond Handoff Datetime" (col H) and
# "Correspond Handback DateTime" (col K) for a9c5f272 (row 3)
$wsZhCn.Range("H3").Value = "2016-08-17 04:41:31"
$wsZhCn.Range("K3").Value = "2016-08-17 04:41:48"

# de-de sheet: "Correspond Handoff Datetime" (col H) and
# "Correspond Handback DateTime" (col K) for a9c5f272 (row 3)
$wsDeDe.Range("H3").Value = "2016-08-17 04:41:37"
$wsDeDe.Range("K3").Value = "2016-08-17 04:41:55"
